$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set H column (Absent) to 1 for rows 3-9, 13-16, and 18
$absentRows = @(3,4,5,6,7,8,9,13,14,15,16,18)
foreach ($r in $absentRows) {
    $ws.Cells.Item($r, 8).Value = 1
}

# Set D (Total Attendance Count) and E (Real) columns to 1 for rows 10, 11, 12, 17
$dualRows = @(10,11,12,17)
foreach ($r in $dualRows) {
    $ws.Cells.Item($r, 4).Value = 1
    $ws.Cells.Item($r, 5).Value = 1
}
